$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the time number format + left/vcenter alignment to the block that will
# hold the new entries (B51:C57) - this reuses the existing "h:mm, no border"
# style already present in the workbook (used by rows 43-50).
$times = $ws.Range("B51:C57")
$times.NumberFormat = "h:mm"
$times.HorizontalAlignment = -4131
$times.VerticalAlignment = -4108

# New work-log rows (Clock GUI / MVVM work). Shared-string entries must be
# created in this exact order to match the authoring order of the edit.
$ws.Range("B51").Value = 0.55208333333333337
$ws.Range("C51").Value = 0.61111111111111105
$ws.Range("D51").Value = "Made basic sprites for GUI in AseSprite."

$ws.Range("B52").Value = 0.61111111111111105
$ws.Range("C52").Value = 0.625
$ws.Range("D52").Value = "Started to plan out Gui architechture."

$ws.Range("B54").Value = 0.72222222222222221
$ws.Range("C54").Value = 0.82986111111111116
$ws.Range("D54").Value = "Marks half finished, Background finished."

$ws.Range("B53").Value = 0.63541666666666663
$ws.Range("C53").Value = 0.6875
$ws.Range("D53").Value = "Working on implementing MVVM for Clock GUI."

$ws.Range("B55").Value = 0.82986111111111116
$ws.Range("C55").Value = 0.85416666666666663
$ws.Range("D55").Value = "ClockGui finished, except the markings."

$ws.Range("B56").Value = 0.89583333333333337
$ws.Range("C56").Value = 0.90972222222222221
$ws.Range("D56").Value = "Creating class to test the visuals."

$ws.Range("B57").Value = 0.90972222222222221
$ws.Range("C57").Value = 0.53125
$ws.Range("D57").Value = "Finished debuging…"

# Update selection to reflect where editing ended up.
[void]$ws.Range("D48").Select()
